# Update the cached "datetimeFigureOut" field text (2023-11-04 -> 2023-11-10)
# on every Date placeholder across the slide master, every slide layout, and
# the handout master, and bump the Title placeholder's default font size on
# the "Title and Content" layout (32pt -> 36pt).

$p = $ppt.ActivePresentation

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "2023-11-04") {
                $tr.Text = "2023-11-10"
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# Every slide layout off the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateShape $layout.Shapes
}

# Handout master
Update-DateShape $p.HandoutMaster.Shapes

# "Title and Content" layout - bump the Title placeholder default size
$titleContentLayout = $master.CustomLayouts.Item(2)
for ($i = 1; $i -le $titleContentLayout.Shapes.Count; $i++) {
    $shp = $titleContentLayout.Shapes.Item($i)
    if ($shp.Name -eq "Title 6") {
        $shp.TextFrame.TextRange.Font.Size = 36
    }
}
